$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates: weekly volume number and date range
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# Transitions (style/type fixups via copy)
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4104) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4104) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4104) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4104) | Out-Null
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4104) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4104) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4104) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4104) | Out-Null
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4104) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4104) | Out-Null

# Final value assignments
$ws.Range("N14").Value = -58.064516129032
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 42
$ws.Range("J15").Value = 38
$ws.Range("K15").Value = 10.526315789473
$ws.Range("M15").Value = 13.513513513513
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 42
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = -10.638297872340
$ws.Range("I16").Value = 427
$ws.Range("J16").Value = 343
$ws.Range("K16").Value = 24.489795918367
$ws.Range("L16").Value = 50.883392226148
$ws.Range("M16").Value = 4.400977995110
$ws.Range("N16").Value = -67.001545595054
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -18.75
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = -18.333333333333
$ws.Range("I17").Value = 677
$ws.Range("J17").Value = 677
$ws.Range("L17").Value = 7.290015847860
$ws.Range("M17").Value = 67.574257425742
$ws.Range("N17").Value = -15.796019900497
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 125
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 68.75
$ws.Range("I18").Value = 279
$ws.Range("J18").Value = 227
$ws.Range("K18").Value = 22.907488986784
$ws.Range("L18").Value = 26.244343891402
$ws.Range("M18").Value = -20.285714285714
$ws.Range("N18").Value = -84.075342465753
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 7.407407407407
$ws.Range("I19").Value = 716
$ws.Range("J19").Value = 497
$ws.Range("K19").Value = 44.064386317907
$ws.Range("L19").Value = 78.553615960099
$ws.Range("M19").Value = 146.896551724138
$ws.Range("N19").Value = 40.117416829745
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 42
$ws.Range("G20").Value = 63
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 413
$ws.Range("J20").Value = 513
$ws.Range("K20").Value = -19.493177387914
$ws.Range("L20").Value = 65.863453815261
$ws.Range("M20").Value = 25.914634146341
$ws.Range("N20").Value = -72.757255936675
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -14.754098360655
$ws.Range("F21").Value = 220
$ws.Range("G21").Value = 242
$ws.Range("H21").Value = -9.090909090909
$ws.Range("I21").Value = 2567
$ws.Range("J21").Value = 2309
$ws.Range("K21").Value = 11.173668254655
$ws.Range("L21").Value = 39.662676822633
$ws.Range("M21").Value = 39.510869565217
$ws.Range("N21").Value = -57.059217129474
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("I22").Value = 37
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 68.181818181818
$ws.Range("L22").Value = 76.190476190476
$ws.Range("M22").Value = 37.037037037037
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -61.538461538461
$ws.Range("I23").Value = 102
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = -7.272727272727
$ws.Range("L23").Value = 7.368421052631
$ws.Range("M23").Value = 47.826086956521
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -4.347826086956
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 24.137931034482
$ws.Range("I24").Value = 1369
$ws.Range("J24").Value = 1023
$ws.Range("K24").Value = 33.822091886608
$ws.Range("L24").Value = 48.481561822125
$ws.Range("M24").Value = 90.668523676880
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -52
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 82
$ws.Range("H25").Value = -32.926829268292
$ws.Range("I25").Value = 811
$ws.Range("J25").Value = 871
$ws.Range("K25").Value = -6.888633754305
$ws.Range("L25").Value = -14.361140443505
$ws.Range("M25").Value = -10.682819383259
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 71
$ws.Range("J26").Value = 55
$ws.Range("K26").Value = 29.090909090909
$ws.Range("L26").Value = 24.561403508771
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -71.428571428571
$ws.Range("J27").Value = 78
$ws.Range("K27").Value = -15.384615384615
$ws.Range("L27").Value = 3.125
$ws.Range("C28").Value = 1
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = -72.727272727272
$ws.Range("I28").Value = 43
$ws.Range("K28").Value = -39.436619718309
$ws.Range("L28").Value = -6.521739130434
$ws.Range("M28").Value = -31.746031746031
$ws.Range("N28").Value = -65.6
$ws.Range("C29").Value = 1
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = -57.142857142857
$ws.Range("I29").Value = 40
$ws.Range("K29").Value = -31.034482758620
$ws.Range("L29").Value = -2.439024390243
$ws.Range("M29").Value = -21.568627450980
$ws.Range("N29").Value = -66.101694915254
